$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new matrix_profile_* columns
$ws.Range("K1").Value = "matrix_profile_identified"
$ws.Range("L1").Value = "matrix_profile_Overlap_merlin"
$ws.Range("M1").Value = "matrix_profilebest_param"
$ws.Range("N1").Value = "matrix_profiletime_taken"

# Match the look of the existing header cells (bold, centered, thin border)
$hdr = $ws.Range("K1:N1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.LineStyle = 1         # xlContinuous

# Row 2 - stdb_308_1.txt
$nums2 = 0..5399
$ws.Range("K2").Value = "[" + ($nums2 -join ", ") + "]"
$ws.Range("L2").Value = 0.06897740989825832
$ws.Range("M2").Value = "{'nbr_of_discord': 2}"
$ws.Range("N2").Value = 0.04688363010063767

# Row 3 - xmitdb_x108_1.txt
$ws.Range("K3").Value = "[4019, 4119, 4325]"
$ws.Range("L3").Value = 0.5542857142857143
$ws.Range("M3").Value = "{'nbr_of_discord': 3}"
$ws.Range("N3").Value = 0.04657184798270464

# Row 4 - mitdb__100_180_1.txt
$ws.Range("K4").Value = "[1794, 1998, 2115]"
$ws.Range("L4").Value = 0.5467224546722456
$ws.Range("M4").Value = "{'nbr_of_discord': 3}"
$ws.Range("N4").Value = 0.0465080130379647

# Row 5 - chfdb_chf01_275_1.txt
$ws.Range("K5").Value = "[2237, 2358, 2472]"
$ws.Range("L5").Value = 0.6299212598425197
$ws.Range("M5").Value = "{'nbr_of_discord': 3}"
$ws.Range("N5").Value = 0.02538140001706779

# Row 6 - ltstdb_20221_43_1.txt
$ws.Range("K6").Value = "[633]"
$ws.Range("L6").Value = 0.415
$ws.Range("M6").Value = "{'nbr_of_discord': 1}"
$ws.Range("N6").Value = 0.02513470803387463

# Row 7 - mitdbx_mitdbx_108.txt
$ws.Range("K7").Value = "[4001, 10352, 10870]"
$ws.Range("L7").Value = 0.8214894826606027
$ws.Range("M7").Value = "{'nbr_of_discord': 3}"
$ws.Range("N7").Value = 0.6312010691035539

# Row 8 - qtdbsele0606.txt
$nums8 = 0..5645
$ws.Range("K8").Value = "[" + ($nums8 -join ", ") + ", "
$ws.Range("L8").Value = 0.005305391604217786
$ws.Range("M8").Value = "{'nbr_of_discord': 1}"
$ws.Range("N8").Value = 0.2957163159735501

# Row 9 - chfdbchf15.txt
$ws.Range("K9").Value = "[2288]"
$ws.Range("L9").Value = 0.8100000000000001
$ws.Range("M9").Value = "{'nbr_of_discord': 1}"
$ws.Range("N9").Value = 0.3001744151115417

# Row 10 - ann_gun_CentroidA_1.txt
$nums10 = 0..5645
$ws.Range("K10").Value = "[" + ($nums10 -join ", ") + ", "
$ws.Range("L10").Value = 0.02597402597402598
$ws.Range("M10").Value = "{'nbr_of_discord': 3}"
$ws.Range("N10").Value = 0.282284809043631

# Row 11 - Patient_respiration.txt
$ws.Range("K11").Value = "[4880, 4968, 5018]"
$ws.Range("L11").Value = 0.45662100456621
$ws.Range("M11").Value = "{'nbr_of_discord': 3}"
$ws.Range("N11").Value = 0.06317739910446107

# Row 12 - dutch_power_demand.txt
$ws.Range("K12").Value = "no"
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = "params"
$ws.Range("N12").Value = 0

# Row 13 - GPS_trajectory_data.csv
$ws.Range("K13").Value = "no"
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = "params"
$ws.Range("N13").Value = 0
